$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains text-formatted numeric-looking strings (e.g. "1.00", "66.338.95").
# Force the cell format to Text before assignment so Excel does not auto-convert these to numbers,
# which preserves the exact original text representation (leading/trailing zeros, dot-grouping, etc).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.338.95'
$ws.Range("E2").Value = '  +5.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.003.61'
$ws.Range("E3").Value = '  +2.79%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.67'
$ws.Range("E5").Value = '  +2.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.80'
$ws.Range("E6").Value = '  +11.84%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +3.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.998.23'
$ws.Range("E9").Value = '  +2.64%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.61'
$ws.Range("E10").Value = '  -4.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.156'
$ws.Range("E11").Value = '  +2.89%  '

$ws.Range("E12").Value = '  +4.62%  '

$ws.Range("E13").Value = '  +5.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.76'
$ws.Range("E14").Value = '  +5.66%  '

$ws.Range("E15").Value = '  -1.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.268.35'
$ws.Range("E16").Value = '  +6.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.501.09'
$ws.Range("E17").Value = '  +3.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.93'
$ws.Range("E18").Value = '  +4.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.004.20'
$ws.Range("E19").Value = '  +3.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '456.39'
$ws.Range("E20").Value = '  +5.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.87'
$ws.Range("E21").Value = '  +5.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  +3.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.35'
$ws.Range("E23").Value = '  +6.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.30'
$ws.Range("E24").Value = '  +4.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("E25").Value = '  +13.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.33'
$ws.Range("E26").Value = '  +2.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.52'
$ws.Range("E27").Value = '  +4.69%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.22'
$ws.Range("E29").Value = '  +17.09%  '

$ws.Range("E30").Value = '  +18.65%  '

$ws.Range("E31").Value = '  -6.12%  '

$ws.Range("E32").Value = '  +4.66%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.22'
$ws.Range("E33").Value = '  +5.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").Value = '  +3.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  +3.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("E37").Value = '  +7.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.17'
$ws.Range("E38").Value = '  +13.61%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.89'
$ws.Range("E39").Value = '  +1.72%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.00'
$ws.Range("E40").Value = '  +0.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.312'
$ws.Range("E41").Value = '  +15.68%  '

$ws.Range("E42").Value = '  +6.97%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '44.03'
$ws.Range("E43").Value = '  +6.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.43'
$ws.Range("E44").Value = '  +3.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '400.57'
$ws.Range("E45").Value = '  +14.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0361'
$ws.Range("E46").Value = '  +6.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.795.28'
$ws.Range("E47").Value = '  +3.20%  '

$ws.Range("E48").Value = '  +1.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.00'
$ws.Range("E50").Value = '  +10.95%  '

$ws.Range("E51").Value = '  +4.26%  '
